# Apply crypto price/volume updates per the commit diff.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "42.677.67"
$ws.Range("E2").Value = "  +0.64%  "
$ws.Range("D3").Value = "2.279.15"
$ws.Range("E3").Value = "  -0.16%  "
$ws.Range("E4").Value = "  +0.09%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "309.86"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -4.17%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "102.77"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +0.39%  "
$ws.Range("E7").Value = "  -1.58%  "
$ws.Range("E9").Value = "  -1.44%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "38.59"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -2.85%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.0894"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  -0.81%  "
$ws.Range("E12").Value = "  -1.39%  "
$ws.Range("E13").Value = "  +0.77%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "0.961"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -0.33%  "
$ws.Range("E15").Value = "  -0.46%  "
$ws.Range("D16").Value = "2.624.73"
$ws.Range("E16").Value = "  -0.11%  "
$ws.Range("D17").Value = "2.271.55"
$ws.Range("E17").Value = "  -0.41%  "
$ws.Range("D18").Value = "42.290.44"
$ws.Range("E18").Value = "  +0.14%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "7.23"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -1.96%  "
$ws.Range("E20").Value = "  -1.50%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "12.97"
$ws.Range("D21").Style = "Normal"
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "72.60"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -0.62%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "3.38"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -7.72%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "261.36"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -2.68%  "
$ws.Range("E25").Value = "  -2.53%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "1.01"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +0.29%  "
$ws.Range("E27").Value = "  -2.25%  "
$ws.Range("B28").Value = "Filecoin"
$ws.Range("C28").Value = "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "6.86"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +12.78%  "
$ws.Range("B29").Value = "Toncoin"
$ws.Range("C29").Value = "https://coinranking.com/coin/67YlI0K1b+toncoin-ton"
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "2.28"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -0.38%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "22.02"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -1.76%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "35.74"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -6.05%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "164.05"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -0.05%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.0846"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -3.56%  "
$ws.Range("E34").Value = "  -3.21%  "
$ws.Range("E35").Value = "  +0.12%  "
$ws.Range("E36").Value = "  -3.83%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "4.46"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -2.87%  "
$ws.Range("E38").Value = "  -2.37%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "3.64"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -1.63%  "
$ws.Range("E40").Value = "  -3.04%  "
$ws.Range("E41").Value = "  +1.92%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "97.34"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +7.48%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.999"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -0.11%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "68.14"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -1.05%  "
$ws.Range("E45").Value = "  -0.72%  "
$ws.Range("E46").Value = "  -3.01%  "
$ws.Range("D47").Value = "1.700.67"
$ws.Range("E47").Value = "  +6.95%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "109.49"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -3.00%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "76.43"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -4.48%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "8.59"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -3.67%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "5.08"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -2.65%  "
